$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A41").Value = "9505052351086"
$ws.Range("C41").Value = "ID Used"
$ws.Range("A42").Value = "9105050675188"
$ws.Range("C42").Value = "ID Used"
$ws.Range("A43").Value = "0105050484184"
$ws.Range("C43").Value = "ID Used"
$ws.Range("A44").Value = "9805054281087"
$ws.Range("C44").Value = "ID Used"
$ws.Range("A45").Value = "9505053023189"
$ws.Range("C45").Value = "ID Used"
$ws.Range("A46").Value = "0105052534085"
$ws.Range("C46").Value = "ID Used"
$ws.Range("A47").Value = "9805050350183"
$ws.Range("C47").Value = "ID Used"
$ws.Range("A48").Value = "9505053406186"
$ws.Range("C48").Value = "ID Used"
$ws.Range("A49").Value = "0105054378085"
$ws.Range("C49").Value = "ID Used"
$ws.Range("A50").Value = "9805052525188"
$ws.Range("C50").Value = "ID Used"
$ws.Range("A51").Value = "9505052006185"
$ws.Range("C51").Value = "ID Used"
$ws.Range("A52").Value = "0105050624185"
$ws.Range("C52").Value = "ID Used"
$ws.Range("A53").Value = "9805051782186"
$ws.Range("C53").Value = "ID Used"
$ws.Range("A54").Value = "9505051784089"
$ws.Range("C54").Value = "ID Used"
$ws.Range("A55").Value = "0105051766084"
$ws.Range("C55").Value = "ID Used"
$ws.Range("A56").Value = "9805054945186"
$ws.Range("C56").Value = "ID Used"
$ws.Range("A57").Value = "9505050169183"
$ws.Range("C57").Value = "ID Used"
$ws.Range("A58").Value = "0105050261186"
$ws.Range("C58").Value = "ID Used"
$ws.Range("A59").Value = "9805051382185"
$ws.Range("C59").Value = "ID Used"
$ws.Range("A60").Value = "9505052228086"
$ws.Range("C60").Value = "ID Used"
$ws.Range("A61").Value = "9105051715181"
$ws.Range("C61").Value = "ID Used"
$ws.Range("A62").Value = "0105051775085"
$ws.Range("C62").Value = "ID Used"
$ws.Range("A63").Value = "98050509840810"
$ws.Range("C63").Value = "ID Used"
$ws.Range("A64").Value = "9505050884187"
$ws.Range("C64").Value = "ID Used"
$ws.Range("A65").Value = "0105051542089"
$ws.Range("C65").Value = "ID Used"
$ws.Range("A66").Value = "9805051204082"
$ws.Range("C66").Value = "ID Used"
$ws.Range("A67").Value = "9505052650081"
$ws.Range("C67").Value = "ID Used"
$ws.Range("A68").Value = "0105054410185"
$ws.Range("C68").Value = "ID Used"
$ws.Range("A69").Value = "9805052254086"
$ws.Range("C69").Value = "ID Used"
$ws.Range("A70").Value = "9505053769088"
$ws.Range("C70").Value = "ID Used"
$ws.Range("A71").Value = "0105050666186"
$ws.Range("C71").Value = "ID Used"
$ws.Range("A72").Value = "9805050541088"
$ws.Range("C72").Value = "ID Used"
$ws.Range("A73").Value = "9505053958087"
$ws.Range("C73").Value = "ID Used"
$ws.Range("A74").Value = "0105050939088"
$ws.Range("C74").Value = "ID Used"
$ws.Range("A75").Value = "9805050872186"
$ws.Range("C75").Value = "ID Used"
$ws.Range("A76").Value = "9505052887188"
$ws.Range("C76").Value = "ID Used"
$ws.Range("A77").Value = "0105052266084"
$ws.Range("C77").Value = "ID Used"
